$wb = $excel.ActiveWorkbook

# Rename sheets
$wb.Worksheets.Item("cr_proportions_total").Name = "proportions_total"
$wb.Worksheets.Item("cr_proportions_by_year").Name = "proportions_year"
$wb.Worksheets.Item("cr_proportions_by_journal").Name = "proportions_journal"
$wb.Worksheets.Item("cr_proportions_by_sample_source").Name = "proportions_sample_source"
$wb.Worksheets.Item("cr_proportions_by_sample_method").Name = "proportions_sample_method"
$wb.Worksheets.Item("cr_proportions_by_sample_plat").Name = "proportions_platform"
$wb.Worksheets.Item("cr_proportions_by_method").Name = "proportions_cr_method"
$wb.Worksheets.Item("cr_proportions_by_method_type").Name = "proportions_cr_type"

# Update header cells on each sheet to match new naming convention
$wb.Worksheets.Item("proportions_total").Range("D1").Value = "proportions_total"
$wb.Worksheets.Item("proportions_year").Range("E1").Value = "proportions_year"
$wb.Worksheets.Item("proportions_journal").Range("F1").Value = "proportions_journal"
$wb.Worksheets.Item("proportions_sample_source").Range("F1").Value = "proportions_sample_source"
$wb.Worksheets.Item("proportions_sample_method").Range("F1").Value = "proportions_sample_method"
$wb.Worksheets.Item("proportions_platform").Range("F1").Value = "proportions_platform"
$wb.Worksheets.Item("proportions_cr_method").Range("G1").Value = "proportions_cr_method"
$wb.Worksheets.Item("proportions_cr_type").Range("H1").Value = "proportions_cr_type"
